# Generate Report for Handback
# Refreshes the handback-status report with newly generated file identifiers,
# hashes and timestamps for the two tracked markdown files.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New data values (replacing the previous run's generated names/timestamps)
# ---------------------------------------------------------------------------
$oldFile1 = "c6644389-fa82-4909-a5dd-6302e36807e0.md"
$oldFile2 = "eda8aa28-9073-4fd6-8907-f1bb0d8dfa47.md"

$newFile1 = "cf0d3d15-07e2-4798-9d1c-82fec8a93031.md"
$newFile2 = "ffff187ec8cc-55eb-458d-bc03-d819daf81daa.md"

$newFile1Path = "e2e\" + $newFile1
$newFile2Path = "e2e\" + $newFile2

$overviewDate = "2016-08-28 17:02:38"

$zhcnXlf = "cf0d3d15-07e2-4798-9d1c-82fec8a93031.ae237857e9cd159d190a9ce1b72edf3cc271d7e0.zh-cn.xlf"
$dedeXlf = "cf0d3d15-07e2-4798-9d1c-82fec8a93031.ae237857e9cd159d190a9ce1b72edf3cc271d7e0.de-de.xlf"

$zhcnHandoffDate  = "2016-08-28 17:02:34"
$zhcnHandbackDate = "2016-08-28 17:02:51"
$dedeHandbackDate = "2016-08-28 17:02:57"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFile1
$wsOverview.Range("B2").Value = $newFile1Path
$wsOverview.Range("G2").Value = $overviewDate

$wsOverview.Range("A3").Value = $newFile2
$wsOverview.Range("B3").Value = $newFile2Path
$wsOverview.Range("G3").Value = $overviewDate

# Rebuild the hyperlinks so the visible display text matches the new paths.
# The underlying link targets are untouched by this refresh (they still
# point at the previously committed source files), so reuse the exact
# original addresses and only change the display text.
$linkOv2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c2d334e1346a701c302cd9e2d1a2ac23166fe964/e2e/" + $oldFile1
$linkOv3 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c2d334e1346a701c302cd9e2d1a2ac23166fe964/e2e/" + $oldFile2

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $linkOv2, "", "", $newFile1Path)
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $linkOv3, "", "", $newFile2Path)

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = $newFile1
$wsZhCn.Range("G2").Value = $zhcnXlf
$wsZhCn.Range("H2").Value = $zhcnHandoffDate
$wsZhCn.Range("I2").Value = $newFile1
$wsZhCn.Range("J2").Value = $zhcnXlf
$wsZhCn.Range("K2").Value = $zhcnHandbackDate

$wsZhCn.Range("A3").Value = $newFile2
$wsZhCn.Range("G3").Value = $zhcnXlf
$wsZhCn.Range("H3").Value = $zhcnHandoffDate
$wsZhCn.Range("I3").Value = $newFile2
$wsZhCn.Range("J3").Value = $zhcnXlf
$wsZhCn.Range("K3").Value = $zhcnHandbackDate

$linkZh2a = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c2d334e1346a701c302cd9e2d1a2ac23166fe964/e2e/" + $oldFile1
$linkZh2b = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/5eaea07cd9b5170427e6c8a9152feb5fd1256280/e2e/" + $oldFile1
$linkZh3a = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c2d334e1346a701c302cd9e2d1a2ac23166fe964/e2e/" + $oldFile2
$linkZh3b = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/5eaea07cd9b5170427e6c8a9152feb5fd1256280/e2e/" + $oldFile2

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $linkZh2a, "", "", $newFile1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $linkZh2b, "", "", $newFile1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $linkZh3a, "", "", $newFile2)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $linkZh3b, "", "", $newFile2)

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = $newFile1
$wsDeDe.Range("G2").Value = $dedeXlf
$wsDeDe.Range("H2").Value = $overviewDate
$wsDeDe.Range("I2").Value = $newFile1
$wsDeDe.Range("J2").Value = $dedeXlf
$wsDeDe.Range("K2").Value = $dedeHandbackDate

$wsDeDe.Range("A3").Value = $newFile2
$wsDeDe.Range("G3").Value = $dedeXlf
$wsDeDe.Range("H3").Value = $overviewDate
$wsDeDe.Range("I3").Value = $newFile2
$wsDeDe.Range("J3").Value = $dedeXlf
$wsDeDe.Range("K3").Value = $dedeHandbackDate

$linkDe2a = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c2d334e1346a701c302cd9e2d1a2ac23166fe964/e2e/" + $oldFile1
$linkDe2b = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/bd90d8482e3bc52fc4ad8764e418a6fe711caaf0/e2e/" + $oldFile1
$linkDe3a = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c2d334e1346a701c302cd9e2d1a2ac23166fe964/e2e/" + $oldFile2
$linkDe3b = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/bd90d8482e3bc52fc4ad8764e418a6fe711caaf0/e2e/" + $oldFile2

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $linkDe2a, "", "", $newFile1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $linkDe2b, "", "", $newFile1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $linkDe3a, "", "", $newFile2)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $linkDe3b, "", "", $newFile2)
